$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.360.67'
$ws.Range("E2").Value = '  -2.37%  '

$ws.Range("D3").Value = '2.646.92'
$ws.Range("E3").Value = '  -3.32%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.96'
$ws.Range("E5").Value = '  -2.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.34'
$ws.Range("E6").Value = '  -4.30%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -1.65%  '

$ws.Range("E9").Value = '  -4.82%  '

$ws.Range("E10").Value = '  -4.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.48'
$ws.Range("E11").Value = '  -3.62%  '

$ws.Range("E12").Value = '  -5.10%  '

$ws.Range("D13").Value = '3.115.42'
$ws.Range("E13").Value = '  -3.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.71'
$ws.Range("E14").Value = '  -4.91%  '

$ws.Range("D15").Value = '62.260.97'
$ws.Range("E15").Value = '  -2.28%  '

$ws.Range("E16").Value = '  -4.26%  '

$ws.Range("D17").Value = '2.653.30'
$ws.Range("E17").Value = '  -3.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.67'
$ws.Range("E18").Value = '  -6.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.54'
$ws.Range("E19").Value = '  -4.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.36'
$ws.Range("E20").Value = '  -4.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.09'
$ws.Range("E21").Value = '  -8.16%  '

$ws.Range("E22").Value = '  -0.12%  '

$ws.Range("E23").Value = '  -3.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.65'
$ws.Range("E24").Value = '  -3.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.00'
$ws.Range("E27").Value = '  -5.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.36'
$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("E29").Value = '  -8.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -1.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.88'
$ws.Range("E31").Value = '  -5.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.68'
$ws.Range("E32").Value = '  -3.57%  '

$ws.Range("E33").Value = '  +0.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.74'
$ws.Range("E34").Value = '  -3.79%  '

$ws.Range("E35").Value = '  -3.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.18'
$ws.Range("E36").Value = '  -4.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.73'
$ws.Range("E37").Value = '  -4.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '333.47'
$ws.Range("E38").Value = '  -3.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.06'
$ws.Range("E39").Value = '  -4.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.908'
$ws.Range("E40").Value = '  -7.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.09'
$ws.Range("E41").Value = '  -1.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.92'
$ws.Range("E42").Value = '  -4.49%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.38'
$ws.Range("E44").Value = '  -6.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.609'
$ws.Range("E45").Value = '  -3.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.81'
$ws.Range("E46").Value = '  -6.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.98'
$ws.Range("E47").Value = '  -0.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0548'
$ws.Range("E48").Value = '  -6.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.88'
$ws.Range("E49").Value = '  -3.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0958'
$ws.Range("E50").Value = '  -4.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0236'
$ws.Range("E51").Value = '  -6.07%  '
